$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Sheet3" worksheet at the end of the workbook and make it
#    the active sheet (tabOTTR 0.2 "DisjointClasses" template w/ list values)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# ---------------------------------------------------------------------------
# 2. Populate Sheet3 with the new OTTR template content
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value = "#OTTR"
$ws3.Range("B1").Value = "prefix"

$ws3.Range("A2").Value = "ex"
$ws3.Range("B2").Value = "http://example.org#"

$ws3.Range("A3").Value = "#OTTR"
$ws3.Range("B3").Value = "end"

$ws3.Range("A4").Value = "#OTTR"
$ws3.Range("B4").Value = "template"
$ws3.Range("C4").Value = "http://candidate.ottr.xyz/owl/axiom/DisjointClasses"

$ws3.Range("A5").Value = "classes"

$ws3.Range("A6").Value = 1

$ws3.Range("A7").Value = "iri+"
# B7 stays empty but carries a (no-op) style so the cell is materialised,
# matching the "range"-column placeholder used by the other templates.
$ws3.Range("B7").Font.Name = "Arial"

$ws3.Range("A8").Value = "ex:ClassA1|ex:ClassA2|ex:ClassA3"
$ws3.Range("A9").Value = "ex:ClassB1|ex:ClassB2"
$ws3.Range("A10").Value = "ex:ClassC1|ex:ClassC2|ex:ClassC3|ex:ClassC4"

$ws3.Range("A11").Value = "#OTTR"
$ws3.Range("B11").Value = "end"

# ---------------------------------------------------------------------------
# 3. Hyperlinks on Sheet3 (mirrors Sheet1 / Sheet1_2: blue, non-underlined,
#    matching the pre-existing "Arial 10 FF0000FF" look instead of Excel's
#    auto "Hyperlink" style)
# ---------------------------------------------------------------------------
$ws3.Hyperlinks.Add($ws3.Range("B2"), "http://example.org/", "", "", "http://example.org#")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "http://candidate.ottr.xyz/owl/axiom/DisjointClasses", "", "", "http://candidate.ottr.xyz/owl/axiom/DisjointClasses")
foreach ($addr in @("B2", "C4")) {
    $r = $ws3.Range($addr)
    $r.Font.Name = "Arial"
    $r.Font.Size = 10
    $r.Font.Color = 16711680
    $r.Font.Underline = $false
}

$ws3.Columns.Item(1).ColumnWidth = 10.6868

# ---------------------------------------------------------------------------
# 4. Fix up the hyperlink display text on Sheet1 / Sheet1_2 (the trailing
#    "#" was missing from the displayed text before)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "http://example.org/", "", "", "http://example.org#")
$ws1.Hyperlinks.Add($ws1.Range("C4"), "http://candidate.ottr.xyz/owl/axiom/SuperObjectMaxCardinality", "", "", "http://candidate.ottr.xyz/owl/axiom/SuperObjectMaxCardinality")
foreach ($addr in @("B2", "C4")) {
    $r = $ws1.Range($addr)
    $r.Font.Name = "Arial"
    $r.Font.Size = 10
    $r.Font.Color = 16711680
    $r.Font.Underline = $false
}

$ws2 = $wb.Worksheets.Item("Sheet1_2")
$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("B2"), "http://example.org/", "", "", "http://example.org#")
$ws2.Hyperlinks.Add($ws2.Range("C7"), "http://candidate.ottr.xyz/owl/axiom/SubObjectMinCardinality", "", "", "http://candidate.ottr.xyz/owl/axiom/SubObjectMinCardinality")
$ws2.Hyperlinks.Add($ws2.Range("B23"), "http://easdfasdfxample/", "", "", "http://easdfasdfxample.org#")
foreach ($addr in @("B2", "C7", "B23")) {
    $r = $ws2.Range($addr)
    $r.Font.Name = "Arial"
    $r.Font.Size = 10
    $r.Font.Color = 16711680
    $r.Font.Underline = $false
}

# ---------------------------------------------------------------------------
# 5. Column widths: columns A & B are now the same width, column C slightly
#    narrower (tabOTTR 0.2 re-layout)
# ---------------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2)) {
    $ws.Columns.Item(1).ColumnWidth = 10.5085
    $ws.Columns.Item(2).ColumnWidth = 10.5085
    $ws.Columns.Item(3).ColumnWidth = 15.6368
}

# ---------------------------------------------------------------------------
# 6. Selections: every sheet now just selects A1, and Sheet3 is the active
#    (front-most) tab.
# ---------------------------------------------------------------------------
$ws1.Range("A1").Select()
$ws2.Range("A1").Select()
$ws3.Range("A10").Select()
$ws3.Activate()

Write-Host "edit complete"
